$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Match Planning")

# Row 4: "take note"
$ws.Range("A4").Value = "take note"
$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B4").Value = 45323

# Row 5: "view notes"
$ws.Range("A5").Value = "view notes"
$ws.Range("B3").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B5").Value = 45323

$ws.Application.CutCopyMode = $false

$ws.Range("B5").Select()
